$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Feature Engineering" in the Machine Learning section,
# right before "Decision Trees" (currently row 20).
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "Machine Learning"
$ws.Range("B20").Value = "Feature Engineering"
$ws.Range("C20").Value = "TBD"

# Insert a new row for "Spectral Decomposition" in the Mathematics section,
# right before "Signal Processing" (now at row 29 after the shift above).
$ws.Rows.Item(29).Insert()
$ws.Range("A29").Value = "Mathematics"
$ws.Range("B29").Value = "Spectral Decomposition"
$ws.Range("C29").Value = "TBD"

# Update the selected cell to match the saved view state.
$ws.Range("B27").Select()
